# Update Financials - refresh scraped yearly data values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83 - Research Development: values shifted one year left, newest year (J) now "NA"
$ws.Range("D83").Value = 1300
$ws.Range("E83").Value = "NA"
$ws.Range("F83").Value = "NA"
$ws.Range("G83").Value = "NA"
$ws.Range("H83").Value = "NA"
$ws.Range("I83").Value = 100
$ws.Range("J83").Value = "NA"

# Row 89 - Total Other Income/Expenses Net
$ws.Range("D89").Value = -8400
$ws.Range("E89").Value = -2600
$ws.Range("F89").Value = -28400
$ws.Range("G89").Value = -23500
$ws.Range("H89").Value = -14500
$ws.Range("I89").Value = -5400
$ws.Range("J89").Value = -6400

# Row 91 - Minority Interest
$ws.Range("D91").Value = -100
$ws.Range("E91").Value = 0
$ws.Range("F91").Value = -200
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = -200
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = -100

# Row 94 - Discontinued Operations
$ws.Range("D94").Value = 3500
$ws.Range("E94").Value = -100
$ws.Range("F94").Value = -2700
$ws.Range("G94").Value = -1600
$ws.Range("H94").Value = -19500
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = "NA"

# Row 100 - Net Income From Continuing Ops
$ws.Range("D100").Value = -2100
$ws.Range("E100").Value = 11600
$ws.Range("F100").Value = 13500
$ws.Range("G100").Value = 1600
$ws.Range("H100").Value = 85700
$ws.Range("I100").Value = 1200
$ws.Range("J100").Value = "NA"

# Row 101 - Net Income Applicable To Common Shares (adjustment row)
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 200
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = "NA"

# Row 102 - Net Income
$ws.Range("D102").Value = -7000
$ws.Range("E102").Value = 8900
$ws.Range("F102").Value = -17600
$ws.Range("G102").Value = -23400
$ws.Range("H102").Value = 51700
$ws.Range("I102").Value = -4200
$ws.Range("J102").Value = 4700
